{"js": "// \"Wireframes version 2.\" -> \"Wireframes version 1.\"\n// The paragraph text is split as \"Versi\" + \"on\" + \" 2\" + \".\" across runs;\n// merge the \"Versi\"/\"on\" pair back into a single \"Version\" run and swap\n// the version digit from 2 to 1 (the trailing \".\" ends up absorbed into\n// the \" 1.\" run, same as the \" 2\" run previously absorbed the separate\n// \".\" run).\nconst body = context.document.body;\n\n// Re-write \"Version\" in place: a same-text replace collapses the\n// \"Versi\" + \"on\" run pair into a single run.\nconst versionHits = body.search(\"Version\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\nif (versionHits.items.length > 0) {\n  versionHits.items[0].insertText(\"Version\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Swap the version number: \" 2.\" -> \" 1.\"\nconst numberHits = body.search(\" 2.\", { matchCase: true });\nnumberHits.load(\"items\");\nawait context.sync();\nif (numberHits.items.length > 0) {\n  numberHits.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# \"Wireframes version 2.\" -> \"Wireframes version 1.\"\n#\n# The paragraph text \"Version 2.\" is split across runs as:\n#   \"Versi\" | \"on\" | (spell-check markers) | \" 2\" | (bookmark _GoBack) | \".\"\n# The target rewrites it to \"Version 1.\" split as:\n#   \"Version\" | (spell-check markers) | \" 1.\" | (bookmark _GoBack)\n# i.e. the \"Versi\"/\"on\" pair collapses into one run, the version digit\n# flips from 2 to 1, and the trailing \".\" run is absorbed into the \" 1\"\n# run ahead of the bookmark (the bookmark itself must stay put).\n\n$d = $word.ActiveDocument\n\n# Step 1: merge the \"Versi\" + \"on\" run pair into a single \"Version\" run.\n# A same-text Find/Replace still re-splits the runs to match the found\n# span, which is exactly the merge we want; this span does not touch the\n# bookmark so it is safe to use Find/Replace here.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Version\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Version\"\n$find1.Execute([ref]\"Version\", [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]\"Version\", [ref]2) | Out-Null\n\n# Step 2: change the version digit and fold the trailing period into the\n# same run (\"2\" -> \"1.\"). Do this as a direct Range.Text assignment (not a\n# Find/Replace) so the _GoBack bookmark sitting right after this text keeps\n# its place instead of being collapsed by the replace.\n$digit = $d.Content\n$fd = $digit.Find\n$fd.ClearFormatting()\n$fd.Text = \"2\"\n$fd.Execute() | Out-Null\n$digit.Text = \"1.\"\n\n# Step 3: the original trailing \".\" (now duplicated, sitting right after\n# the bookmark) is redundant - remove it.\n$paraEnd = $d.Paragraphs(1).Range.End\n$tail = $d.Range($paraEnd - 2, $paraEnd - 1)\nif ($tail.Text -eq \".\") {\n  $tail.Delete()\n}\n"}
